$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E for data rows to remain as Text so that values
# like "1.000", "0.9996", "0.3630", "6.000" etc. keep their exact formatting
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.355.00"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "1.626.02"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "0.9996"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").Value = "302.61"
$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("D7").Value = "0.3763"
$ws.Range("E7").Value = "  +0.80%  "

$ws.Range("D8").Value = "0.3630"
$ws.Range("E8").Value = "  +0.48%  "

$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("D10").Value = "0.08145"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("E11").Value = "  -1.40%  "

$ws.Range("D12").Value = "0.9994"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("E13").Value = "  -2.12%  "

$ws.Range("D14").Value = "6.473"
$ws.Range("E14").Value = "  -1.61%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.309"
$ws.Range("E15").Value = "  +0.60%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.00001238"
$ws.Range("E16").Value = "  -2.18%  "

$ws.Range("D17").Value = "1.618.33"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").Value = "94.21"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").Value = "0.06920"
$ws.Range("E19").Value = "  +0.66%  "

$ws.Range("E20").Value = "  -2.99%  "

$ws.Range("D21").Value = "6.541"
$ws.Range("E21").Value = "  +0.53%  "

$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("E23").Value = "  -1.66%  "

$ws.Range("D24").Value = "23.347.95"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("D25").Value = "2.488"
$ws.Range("E25").Value = "  +3.17%  "

$ws.Range("D26").Value = "3.077"
$ws.Range("E26").Value = "  +2.48%  "

$ws.Range("D27").Value = "21.14"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").Value = "150.15"
$ws.Range("E28").Value = "  -0.98%  "

$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").Value = "132.86"
$ws.Range("E30").Value = "  -1.66%  "

$ws.Range("D31").Value = "1.804.64"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").Value = "6.718"
$ws.Range("E32").Value = "  -0.62%  "

$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "2.130"
$ws.Range("E33").Value = "  -6.19%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.052"
$ws.Range("E34").Value = "  +11.15%  "

$ws.Range("D35").Value = "11.13"
$ws.Range("E35").Value = "  +8.05%  "

$ws.Range("D36").Value = "0.02761"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("D37").Value = "0.08769"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").Value = "0.2486"
$ws.Range("E38").Value = "  -1.37%  "

$ws.Range("D39").Value = "0.07130"
$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("D40").Value = "6.000"
$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("D41").Value = "0.6987"
$ws.Range("E41").Value = "  -0.38%  "

$ws.Range("D42").Value = "1.338"
$ws.Range("E42").Value = "  -1.93%  "

$ws.Range("D43").Value = "15.88"
$ws.Range("E43").Value = "  -1.04%  "

$ws.Range("D44").Value = "12.05"
$ws.Range("E44").Value = "  -2.62%  "

$ws.Range("D45").Value = "0.6459"
$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").Value = "0.9992"

$ws.Range("D47").Value = "2.273"
$ws.Range("E47").Value = "  -1.96%  "

$ws.Range("D48").Value = "3.959"
$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").Value = "127.33"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").Value = "1.194"
$ws.Range("E51").Value = "  -0.03%  "
